$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(3, 4).Value = 0.38436940312385559
$ws.Cells.Item(3, 5).Value = 0.078547939658164978
$ws.Cells.Item(3, 8).Value = 0.23041544854640961
$ws.Cells.Item(3, 9).Value = 0.53832334280014038
$ws.Cells.Item(4, 4).Value = 0.39489689469337463
$ws.Cells.Item(4, 5).Value = 0.078298427164554596
$ws.Cells.Item(4, 8).Value = 0.24143198132514954
$ws.Cells.Item(4, 9).Value = 0.54836183786392212
$ws.Cells.Item(5, 4).Value = 0.41082605719566345
$ws.Cells.Item(5, 5).Value = 0.072752475738525391
$ws.Cells.Item(5, 8).Value = 0.26823121309280396
$ws.Cells.Item(5, 9).Value = 0.55342090129852295
$ws.Cells.Item(7, 4).Value = 0.89445489645004272
$ws.Cells.Item(7, 5).Value = 0.018512062728404999
$ws.Cells.Item(7, 8).Value = 0.85817122459411621
$ws.Cells.Item(7, 9).Value = 0.93073856830596924
$ws.Cells.Item(8, 4).Value = 0.89466220140457153
$ws.Cells.Item(8, 5).Value = 0.01751946285367012
$ws.Cells.Item(8, 8).Value = 0.86032402515411377
$ws.Cells.Item(8, 9).Value = 0.9290003776550293
$ws.Cells.Item(9, 4).Value = 0.89355921745300293
$ws.Cells.Item(9, 5).Value = 0.016960123553872108
$ws.Cells.Item(9, 8).Value = 0.86031734943389893
$ws.Cells.Item(9, 9).Value = 0.92680108547210693
$ws.Cells.Item(11, 4).Value = 0.23652034997940063
$ws.Cells.Item(11, 5).Value = 0.032814193516969681
$ws.Cells.Item(11, 8).Value = 0.17220452427864075
$ws.Cells.Item(11, 9).Value = 0.30083617568016052
$ws.Cells.Item(12, 4).Value = 0.24002310633659363
$ws.Cells.Item(12, 5).Value = 0.030129756778478622
$ws.Cells.Item(12, 8).Value = 0.18096877634525299
$ws.Cells.Item(12, 9).Value = 0.29907742142677307
$ws.Cells.Item(13, 4).Value = 0.24305523931980133
$ws.Cells.Item(13, 5).Value = 0.031533610075712204
$ws.Cells.Item(13, 8).Value = 0.18124936521053314
$ws.Cells.Item(13, 9).Value = 0.30486112833023071
$ws.Cells.Item(15, 4).Value = 0.3303561806678772
$ws.Cells.Item(15, 5).Value = 0.068532422184944153
$ws.Cells.Item(15, 8).Value = 0.19603262841701508
$ws.Cells.Item(15, 9).Value = 0.46467971801757812
$ws.Cells.Item(16, 4).Value = 0.34271252155303955
$ws.Cells.Item(16, 5).Value = 0.066865712404251099
$ws.Cells.Item(16, 8).Value = 0.21165572106838226
$ws.Cells.Item(16, 9).Value = 0.47376930713653564
$ws.Cells.Item(17, 4).Value = 0.35777807235717773
$ws.Cells.Item(17, 5).Value = 0.060750715434551239
$ws.Cells.Item(17, 8).Value = 0.23870666325092316
$ws.Cells.Item(17, 9).Value = 0.47684946656227112
$ws.Cells.Item(19, 4).Value = 0.076747387647628784
$ws.Cells.Item(19, 5).Value = 0.031236883252859116
$ws.Cells.Item(19, 8).Value = 0.015523096546530724
$ws.Cells.Item(19, 9).Value = 0.13797168433666229
$ws.Cells.Item(20, 4).Value = 0.071647018194198608
$ws.Cells.Item(20, 5).Value = 0.032153081148862839
$ws.Cells.Item(20, 8).Value = 0.0086269788444042206
$ws.Cells.Item(20, 9).Value = 0.1346670538187027
$ws.Cells.Item(21, 4).Value = 0.072508752346038818
$ws.Cells.Item(21, 5).Value = 0.032341912388801575
$ws.Cells.Item(21, 8).Value = 0.0091186044737696648
$ws.Cells.Item(21, 9).Value = 0.1358989030122757
